$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("X2").Value = 29
$ws.Range("Z2").Value = 44
$ws.Range("AC2").Value = 11.5
$ws.Range("AD2").Value = 19
$ws.Range("AE2").Value = 46
$ws.Range("AF2").Value = 14.5
$ws.Range("AG2").Value = 10.5
$ws.Range("AH2").Value = 15.5
$ws.Range("AI2").Value = 48
$ws.Range("AK2").Value = 16.5
$ws.Range("AL2").Value = 24
$ws.Range("AM2").Value = 60
$ws.Range("AN2").Value = 6.8
$ws.Range("AO2").Value = 30

# Row 4
$ws.Range("F4").Value = 1.76
$ws.Range("G4").Value = 1.95
$ws.Range("H4").Value = 3.85
$ws.Range("I4").Value = 5.1
$ws.Range("J4").Value = 3.9
$ws.Range("K4").Value = 5.2
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 4.5
$ws.Range("O4").Value = 1.23
$ws.Range("P4").Value = 2.22
$ws.Range("Q4").Value = 1.57
$ws.Range("R4").Value = 1.48
$ws.Range("S4").Value = 2.64
$ws.Range("T4").Value = 1.65
$ws.Range("U4").Value = 2.22
$ws.Range("V4").Value = 1.25
$ws.Range("W4").Value = 2.04
$ws.Range("Y4").Value = 990
$ws.Range("AB4").Value = 1000
$ws.Range("AC4").Value = 1000
$ws.Range("AD4").Value = 1000
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 1000

# Row 6
$ws.Range("F6").Value = 11
$ws.Range("G6").Value = 18
$ws.Range("H6").Value = 1.21
$ws.Range("I6").Value = 1.33
$ws.Range("J6").Value = 5.9
$ws.Range("K6").Value = 8.2

# Row 8
$ws.Range("F8").Value = 2.28
$ws.Range("G8").Value = 2.3
$ws.Range("H8").Value = 3.65
$ws.Range("I8").Value = 3.7
$ws.Range("N8").Value = 3.55
$ws.Range("P8").Value = 1.83
$ws.Range("R8").Value = 1.32
$ws.Range("V8").Value = 1.37
$ws.Range("W8").Value = 1.77
$ws.Range("AF8").Value = 13

# Row 9
$ws.Range("N9").Value = 5.3
$ws.Range("P9").Value = 2.48
$ws.Range("R9").Value = 1.57
$ws.Range("T9").Value = 1.76
$ws.Range("U9").Value = 2.24
$ws.Range("AG9").Value = 9.2

# Row 10
$ws.Range("I10").Value = 4.1
$ws.Range("J10").Value = 3.95
$ws.Range("V10").Value = 1.32

# Row 11
$ws.Range("P11").Value = 1.94
$ws.Range("T11").Value = 1.81

# Row 12
$ws.Range("H12").Value = 1.43
$ws.Range("I12").Value = 1.44
$ws.Range("T12").Value = 1.94
$ws.Range("U12").Value = 2
$ws.Range("V12").Value = 3.25
$ws.Range("Z12").Value = 8.4
$ws.Range("AB12").Value = 32
$ws.Range("AD12").Value = 9.8
$ws.Range("AJ12").Value = 260
$ws.Range("AO12").Value = 6
